# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# (GitHub Actions scheduled scrape update: prices + 1h volume deltas, plus a
#  rank swap between InjectiveProtocol and Aave at rows 40/41.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.421.29'
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = '2.160.06'
$ws.Range("E3").Value = '  +3.49%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''229.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '''0.622'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("D7").Value = '''63.16'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.77%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.396'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("D10").Value = '''0.0868'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '''16.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.83%  '
$ws.Range("D13").Value = '2.481.38'
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").Value = '''22.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '''0.815'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = '''5.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '2.166.79'
$ws.Range("E17").Value = '  +3.65%  '
$ws.Range("D18").Value = '39.392.02'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("D19").Value = '''72.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("D20").Value = '''6.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("D21").Value = '0.0₃0853'
$ws.Range("E21").Value = '  +1.70%  '
$ws.Range("D22").Value = '''228.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").Value = '''2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").Value = '''9.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("D27").Value = '''171.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").Value = '''2.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.32%  '
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").Value = '''4.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.84%  '
$ws.Range("D34").Value = '''4.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("E35").Value = '  +9.26%  '
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("E37").Value = '  +2.24%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''103.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''18.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("E42").Value = '  +1.78%  '
$ws.Range("D43").Value = '1.533.69'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("E44").Value = '  +6.10%  '
$ws.Range("D45").Value = '''0.0932'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.32%  '
$ws.Range("D46").Value = '''1.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.08%  '
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").Value = '''7.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").Value = '''4.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("D50").Value = '2.365.74'
$ws.Range("E50").Value = '  +3.40%  '
$ws.Range("E51").Value = '  +0.23%  '
